{"js": "// Add a new \"CH\u01af\u01a0NG 3: T\u1ed5ng k\u1ebft\" entry to the M\u1ee4C L\u1ee4C (table of contents)\n// list, right after the existing \"CH\u01af\u01a0NG 2: C\u00e1c l\u1ec7nh c\u01a1 b\u1ea3n v\u1ec1 Git\" entry\n// (i.e. as the new last paragraph of the document body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target insertion point is the very last paragraph of the body\n// (the \"CH\u01af\u01a0NG 2 ... Git\" list item, immediately before the section break).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Inserting a sibling paragraph \"After\" the last paragraph naturally\n// inherits that paragraph's formatting (style \"ListParagraph\" and\n// 1.5-line spacing), exactly like pressing Enter at the end of the line\n// in Word, so no extra formatting calls are required.\nconst newParagraph = lastParagraph.insertParagraph(\"CH\u01af\u01a0NG 3: T\u1ed5ng k\u1ebft\", \"After\");\nnewParagraph.load(\"text,style\");\nawait context.sync();\n", "ps1": "# Add a new \"CH\u01af\u01a0NG 3: T\u1ed5ng k\u1ebft\" entry to the M\u1ee4C L\u1ee4C (table of contents)\n# list, right after the existing \"CH\u01af\u01a0NG 2: C\u00e1c l\u1ec7nh c\u01a1 b\u1ea3n v\u1ec1 Git\" entry\n# (i.e. as the new last paragraph of the document body).\n$d = $word.ActiveDocument\n\n# The target insertion point is the very last paragraph of the body\n# (the \"CH\u01af\u01a0NG 2 ... Git\" list item, immediately before the section break).\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n# The freshly split paragraph mark already inherited the \"ListParagraph\"\n# style and 1.5-line spacing from the paragraph it was split off of (just\n# like pressing Enter at the end of the line in Word), so only the text\n# needs to be written.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"CH\u01af\u01a0NG 3: T\u1ed5ng k\u1ebft\"\n"}
